$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.844.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").Value = "1.706.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'311.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").Value = "'0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'0.3749"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.33%  "

$ws.Range("D8").Value = "'49.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.91%  "

$ws.Range("D9").Value = "'0.3448"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").Value = "'1.215"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.87%  "

$ws.Range("D11").Value = "'0.07544"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.40%  "

$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "'21.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.86%  "

$ws.Range("D14").Value = "'6.311"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.42%  "

$ws.Range("D15").Value = "'7.072"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.50%  "

$ws.Range("D16").Value = "1.707.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").Value = "'0.00001131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.66%  "

$ws.Range("D18").Value = "'0.06722"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").Value = "'0.9995"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").Value = "'84.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("D21").Value = "'17.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.91%  "

$ws.Range("D22").Value = "'6.397"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.40%  "

$ws.Range("D23").Value = "'13.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.39%  "

$ws.Range("D24").Value = "24.864.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.82%  "

$ws.Range("D25").Value = "'2.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").Value = "'2.789"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.25%  "

$ws.Range("D27").Value = "'20.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.72%  "

$ws.Range("D28").Value = "'150.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.60%  "

$ws.Range("D29").Value = "'132.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.99%  "

$ws.Range("D30").Value = "1.899.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "

$ws.Range("D31").Value = "'1.241"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +26.39%  "

$ws.Range("D32").Value = "'6.848"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.61%  "

$ws.Range("D33").Value = "'4.239"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.71%  "

$ws.Range("D34").Value = "'13.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.16%  "

$ws.Range("D35").Value = "'0.08820"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.06%  "

$ws.Range("D36").Value = "'1.766"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.86%  "

$ws.Range("D37").Value = "'5.652"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.90%  "

$ws.Range("D38").Value = "'9.292"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.34%  "

$ws.Range("D39").Value = "'0.06658"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("D40").Value = "'0.02409"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "

$ws.Range("D41").Value = "'0.2231"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.18%  "

$ws.Range("D42").Value = "'1.279"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("D43").Value = "'0.6456"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.21%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").Value = "'13.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.19%  "

$ws.Range("D46").Value = "'0.6150"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.90%  "

$ws.Range("D47").Value = "'3.826"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.49%  "

$ws.Range("D48").Value = "'2.131"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.82%  "

$ws.Range("D49").Value = "'129.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.63%  "

$ws.Range("D50").Value = "'0.07298"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").Value = "'80.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.35%  "
